$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" counts for three events
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 135
$wsExhibit.Range("F4").Value = 3587
$wsExhibit.Range("F7").Value = 432

# Sheet "全部类型" (sheet4): same three events, different row for the third one
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 135
$wsAll.Range("F4").Value = 3587
$wsAll.Range("F9").Value = 432
